# Weekly Fruta/Hortaliza update: a new price-reporting row for
# "Femacal de La Calera" / Ciboulette is inserted at row 256, pushing
# the existing rows 256-371 down to 257-372 (dimension grows from
# A1:R371 to A1:R372, matching every row's other columns which repeat
# the same Mercado/Categoria/Calidad/etc. values, date = D, volume = J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 256; Excel shifts rows 256..371 down to 257..372
# and carries the row-256 formatting (e.g. the date style on column D).
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row 256 with the new weekly record.
$ws.Range("A256").Value = 3
$ws.Range("B256").Value = "Femacal de La Calera"
$ws.Range("C256").Value = "Coquimbo"
$ws.Range("D256").Value = 44839
$ws.Range("E256").Value = 5
$ws.Range("F256").Value = 100112039
$ws.Range("G256").Value = "Ciboulette"
$ws.Range("H256").Value = "Sin especificar"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 120
$ws.Range("K256").Value = 1500
$ws.Range("L256").Value = 1500
$ws.Range("M256").Value = 1500
$ws.Range("N256").Value = "$/docena de atados"
$ws.Range("O256").Value = "Provincia de Quillota"
$ws.Range("P256").Value = 500
$ws.Range("Q256").Value = 3
$ws.Range("R256").Value = "Hortaliza"
